# Auto-generated script applying 2024-10-18 YTD violent crime data update
# Updates column K (2024 year-to-date totals) across Citywide Totals,
# By Neighborhood, and per-neighborhood sheets, per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 6339
$ws.Range("K3").Value = 6545
$ws.Range("B4").Value = 1415
$ws.Range("E4").Value = 1641
$ws.Range("K4").Value = 1358
$ws.Range("K5").Value = 464
$ws.Range("K6").Value = 7209
$ws.Range("B7").Value = 18512
$ws.Range("E7").Value = 21012
$ws.Range("K7").Value = 21915

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K6").Value = 155
$ws.Range("K7").Value = 645
$ws.Range("K8").Value = 1436
$ws.Range("K10").Value = 128
$ws.Range("K11").Value = 407
$ws.Range("K19").Value = 640
$ws.Range("K20").Value = 521
$ws.Range("K21").Value = 71
$ws.Range("K24").Value = 68
$ws.Range("K29").Value = 1190
$ws.Range("K31").Value = 244
$ws.Range("K33").Value = 959
$ws.Range("K36").Value = 280
$ws.Range("K37").Value = 743
$ws.Range("K40").Value = 49
$ws.Range("K41").Value = 153
$ws.Range("K42").Value = 810
$ws.Range("K46").Value = 43
$ws.Range("K48").Value = 273
$ws.Range("K50").Value = 103
$ws.Range("K51").Value = 281
$ws.Range("K52").Value = 576
$ws.Range("K53").Value = 281
$ws.Range("K54").Value = 431
$ws.Range("K55").Value = 240
$ws.Range("K60").Value = 129
$ws.Range("B63").Value = 323
$ws.Range("E63").Value = 277
$ws.Range("K63").Value = 58
$ws.Range("K64").Value = 138
$ws.Range("K65").Value = 515
$ws.Range("K67").Value = 858
$ws.Range("K69").Value = 49
$ws.Range("K72").Value = 111
$ws.Range("K76").Value = 299
$ws.Range("K77").Value = 150
$ws.Range("K78").Value = 248
$ws.Range("K80").Value = 77
$ws.Range("K84").Value = 176
$ws.Range("K85").Value = 1019
$ws.Range("K86").Value = 132
$ws.Range("K88").Value = 233
$ws.Range("K89").Value = 323
$ws.Range("K90").Value = 206
$ws.Range("K94").Value = 293
$ws.Range("K95").Value = 360
$ws.Range("K96").Value = 231
$ws.Range("K97").Value = 172
$ws.Range("K98").Value = 110
$ws.Range("K99").Value = 361
$ws.Range("B101").Value = 18512
$ws.Range("E101").Value = 21012
$ws.Range("K101").Value = 21915

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K2").Value = 74
$ws.Range("K4").Value = 13
$ws.Range("K7").Value = 231

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 210
$ws.Range("K3").Value = 212
$ws.Range("K7").Value = 645

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K3").Value = 106
$ws.Range("K7").Value = 407

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K3").Value = 98
$ws.Range("K7").Value = 323

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 334
$ws.Range("K3").Value = 350
$ws.Range("K7").Value = 1019

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 151
$ws.Range("K6").Value = 209
$ws.Range("K7").Value = 576

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 49

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K3").Value = 76
$ws.Range("K6").Value = 118
$ws.Range("K7").Value = 281

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 395
$ws.Range("K3").Value = 436
$ws.Range("K6").Value = 482
$ws.Range("K7").Value = 1436

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 247
$ws.Range("K3").Value = 342
$ws.Range("K5").Value = 26
$ws.Range("K6").Value = 295
$ws.Range("K7").Value = 959

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K6").Value = 82
$ws.Range("K7").Value = 360

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 210
$ws.Range("K7").Value = 743

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 174
$ws.Range("K6").Value = 183
$ws.Range("K7").Value = 515

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 93
$ws.Range("K3").Value = 150
$ws.Range("K7").Value = 361

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 80
$ws.Range("K7").Value = 244

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 235
$ws.Range("K3").Value = 313
$ws.Range("K7").Value = 858

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K3").Value = 69
$ws.Range("K7").Value = 176

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K2").Value = 68
$ws.Range("K7").Value = 431

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 336
$ws.Range("K3").Value = 425
$ws.Range("K6").Value = 343
$ws.Range("K7").Value = 1190

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K3").Value = 65
$ws.Range("K7").Value = 273

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 192
$ws.Range("K6").Value = 209
$ws.Range("K7").Value = 640

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K6").Value = 155
$ws.Range("K7").Value = 299

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K3").Value = 44
$ws.Range("K6").Value = 41
$ws.Range("K7").Value = 155

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K2").Value = 51
$ws.Range("K7").Value = 153

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 219
$ws.Range("K3").Value = 246
$ws.Range("K5").Value = 11
$ws.Range("K7").Value = 810

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K2").Value = 40
$ws.Range("K7").Value = 128

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K3").Value = 62
$ws.Range("K7").Value = 248

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K3").Value = 68
$ws.Range("K7").Value = 240

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("K3").Value = 15
$ws.Range("K7").Value = 68

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 43

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("K6").Value = 42
$ws.Range("K7").Value = 71

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K2").Value = 28
$ws.Range("K6").Value = 52
$ws.Range("K7").Value = 138

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 180
$ws.Range("K7").Value = 521

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K3").Value = 85
$ws.Range("K7").Value = 280

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 131
$ws.Range("K7").Value = 293

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 110

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 103

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K2").Value = 35
$ws.Range("K7").Value = 172

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K2").Value = 60
$ws.Range("K7").Value = 233

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K3").Value = 22
$ws.Range("K7").Value = 132

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K6").Value = 49
$ws.Range("K7").Value = 206

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K6").Value = 94
$ws.Range("K7").Value = 281

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K3").Value = 38
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K2").Value = 24
$ws.Range("K7").Value = 111

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K3").Value = 58
$ws.Range("K7").Value = 150

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("K6").Value = 36
$ws.Range("K7").Value = 77

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("K3").Value = 22
$ws.Range("K7").Value = 49
